$d = $word.ActiveDocument

# Build the OOXML package fragment for the two new paragraphs that were
# appended at the end of the document body (just before <w:sectPr/>):
#   1) a paragraph describing the Laravel-cors package usage
#   2) a trailing empty paragraph
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p>
  <w:r><w:t>Laravel-</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>cors</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> =&gt; </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>packeg</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> used in </w:t></w:r>
  <w:r><w:t xml:space="preserve">project Laravel </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>vuejs</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> separated back-end and front-end</w:t></w:r>
</w:p>
<w:p/>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
"@

# Insert right at the very end of the document's story (after the last
# run of the last existing paragraph, before the sectPr) so the
# pre-existing final paragraph is left completely untouched.
$endPos = $d.Content.End
$insertionPoint = $d.Range($endPos, $endPos)
$insertionPoint.InsertXML($xml)
